# Teens Club Report -> Version 6 update
# - Insert two new columns (P:Q) for "Teen Club - CD4 Drawn" and
#   "Teen Club - Urine LAM Drawn" (pushes old EAC/Actual-attendance columns
#   right by two, from P/Q to R/S, and the trailing helper column from S to U)
# - Populate the new columns with data (mostly 0, with two exceptions)
# - Shift every visit date back by 366 days (one year, Version 6 dataset)
# - Update the active-cell selection

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert the two new columns right before the old "EAC" column (P).
# This automatically shifts P->R, Q->S, S->U, updates dimension/spans/cols
# and carries over per-row styling for the new cells.
$ws.Columns("P:Q").Insert()

# New column headers (new shared strings)
$ws.Range("P1").Value = "Teen Club - CD4 Drawn"
$ws.Range("Q1").Value = "Teen Club - Urine LAM Drawn"

# New column data - defaults to 0 for every data row
$newP = @{2=0; 3=0; 4=0; 5=0; 6=9; 7=0; 8=0; 9=0; 10=0; 11=0; 12=0; 13=0; 14=0; 15=0}
$newQ = @{2=0; 3=0; 4=0; 5=0; 6=0; 7=0; 8=2; 9=0; 10=0; 11=0; 12=0; 13=0; 14=0; 15=0}

foreach ($row in 2..15) {
    $ws.Range("P$row").Value = $newP[$row]
    $ws.Range("Q$row").Value = $newQ[$row]
}

# Shift every visit date back by 366 days (Version 6 Teens Club Report)
$dates = @{2=45272; 3=45262; 4=45265; 5=45269; 6=45269; 7=45277; 8=45277; 9=45262; 10=45269; 11=45263; 12=45270; 13=45262; 14=45270; 15=45269}
foreach ($row in 2..15) {
    $ws.Range("B$row").Value = $dates[$row]
}

# Match the author's final selection
$ws.Range("Q8").Select()
